$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert a new "Collector" column at H (old H=Petri_density shifts to I, etc.) ---
$ws.Columns("H").Insert()

# Header for new column (text-formatted, like the other date/code columns)
$ws.Cells.Item(1, 8).NumberFormat = "@"
$ws.Cells.Item(1, 8).Value2 = "Collector"

# --- Populate Collector column (H) for existing rows 2-5 ---
$ws.Cells.Item(2, 8).NumberFormat = "@"
$ws.Cells.Item(2, 8).Value2 = "AN"
$ws.Cells.Item(3, 8).NumberFormat = "@"
$ws.Cells.Item(3, 8).Value2 = "PE"
$ws.Cells.Item(4, 8).NumberFormat = "@"
$ws.Cells.Item(4, 8).Value2 = "PE"
$ws.Cells.Item(5, 8).NumberFormat = "@"
$ws.Cells.Item(5, 8).Value2 = "AN"
# rows 6-8 have no collector (left blank)

# --- Append new rows 9-15 with the next collection batch (copied layout from rows 2-8) ---
$newRows = @(
    @{ Row=9;  A="Grant";   B="Ferris";     C="Apple"; G=1; H="AN"; I=23 },
    @{ Row=10; A="Grant";   B="Ferris";     C="Apple"; G=2; H="PE"; I=21 },
    @{ Row=11; A="Grant";   B="OG";         C="Apple"; G=1; H="AN"; I=16 },
    @{ Row=12; A="Grant";   B="OG";         C="Apple"; G=2; H="PE"; I=45 },
    @{ Row=13; A="Lansing"; B="MSU";        C="Apple"; G=0; H=$null; I=0 },
    @{ Row=14; A="Lansing"; B="Coll Halls"; C="Haw";   G=0; H=$null; I=0 },
    @{ Row=15; A="Lansing"; B="MilfordST";  C="Haw";   G=1; H="AN"; I=2 }
)

foreach ($rowData in $newRows) {
    $r = $rowData.Row
    $ws.Cells.Item($r, 1).Value2 = $rowData.A
    $ws.Cells.Item($r, 2).Value2 = $rowData.B
    $ws.Cells.Item($r, 3).Value2 = $rowData.C
    $ws.Cells.Item($r, 4).Value2 = "2017-08-24"
    $ws.Cells.Item($r, 6).Value2 = "2017-09-06"
    $ws.Cells.Item($r, 7).Value2 = $rowData.G
    if ($rowData.H) {
        $ws.Cells.Item($r, 8).NumberFormat = "@"
        $ws.Cells.Item($r, 8).Value2 = $rowData.H
    }
    $ws.Cells.Item($r, 9).Value2 = $rowData.I
}

# --- Set the new rows' collection-larvae date (column E) to 2017-08-25 ---
for ($r = 9; $r -le 15; $r++) {
    $ws.Cells.Item($r, 5).Value2 = "2017-08-25"
}

# --- Update Collect_fruit date (column D) for ALL rows (existing + new) ---
for ($r = 2; $r -le 15; $r++) {
    $ws.Cells.Item($r, 4).Value2 = "2017-08-21"
}

# --- Set the new rows' fridge date (column F) to 2017-09-07 ---
for ($r = 9; $r -le 15; $r++) {
    $ws.Cells.Item($r, 6).Value2 = "2017-09-07"
}

# --- Restore selection to J15 as in the edited workbook ---
$ws.Range("J15").Select()
